# Rebuild templates with native PPTX elements (editable text, shapes)
# - remove the background picture
# - turn the old "body_text" placeholder shape into the disclosures paragraph block ("TextBox 1")
# - turn the old "slide_title" placeholder shape into the "Disclosures" heading ("TextBox 2")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- remove the background picture (p:pic, id 2) -----------------------
$s.Shapes.Item(1).Delete()

# After the delete, the two remaining placeholder shapes are:
#   idx1 = "slide_title" (id 3)
#   idx2 = "body_text"   (id 4)
$shpTitle = $s.Shapes.Item(1)
$shpBody  = $s.Shapes.Item(2)

# Put the (future) body textbox first in the shape order, matching the
# target deck where "TextBox 1" (disclosures paragraphs) precedes
# "TextBox 2" (the "Disclosures" heading).
$shpBody.ZOrder(1)   # msoSendToBack

$shpBody  = $s.Shapes.Item(1)
$shpTitle = $s.Shapes.Item(2)

# --- TextBox 1: disclosures paragraph block -----------------------------
$para1 = 'Potomac Fund Management (“Potomac”) is an SEC registered investment adviser located in Bethesda, Maryland. Registration does not imply a certain level of skill or training, nor is it an endorsement by the SEC. This material is for general informational purposes only and does not constitute investment advice, tax advice, or a recommendation regarding any specific product, security, strategy, or investment decision. Readers should not assume that any discussion or information applies to their individual circumstances. This communication does not constitute an offer to buy or sell any security or a solicitation to provide personalized investment advice for compensation. Nothing herein should be construed as individualized or tailored advice delivered over the internet.'
$para2 = 'Opinions expressed are current as of the date of publication and may change without notice. Information obtained from third party sources is believed to be reliable, but Potomac does not guarantee its accuracy or completeness and is not responsible for any third party content referenced or linked in this material.'
$para3 = 'Investing involves risk, including the possible loss of principal. Past performance does not guarantee future results. For additional important disclosures, please visit potomac.com/disclosures.'
$para4 = 'Performance results of Potomac strategies reﬂect the composite performance of all fully discretionary portfolios managed by Potomac according to the strategy subject to policies that may require the exclusion of certain accounts. All returns are time-weighted and reﬂect the reinvestment of dividends and capital gain distributions. Gross performance returns do not reflect the payment of investment advisory fees but reflect the underlying fund management fees, other fund (administrative) expenses, and redemption or 12b1 (fund marketing) fees, if any. Net performance reﬂects the deduction of a model fee (the highest investment advisory fee charged by Potomac), underlying fund management fees, other fund (administrative) expenses and, if any, redemption or 12b1 (fund marketing) fees. Net of fee returns are calculated using a model fee of 2.5%. The model fee, applied monthly, is the highest fee that may be or has been charged to an investor in this composite. Actual investment advisory fees incurred may vary. Past performance does not guarantee future results. There is no guarantee that any investment strategy or account will be proﬁtable or will avoid loss. Individual investors’ objectives, ﬁnancial situations, their speciﬁc instructions, or restrictions on investments, or the time at which an account is opened, or additions are made may result in different trades and returns. Performance for the strategy presented may differ materially (more or less) from the performance of the comparable benchmark and other Potomac investment strategies. Market and economic conditions could change in the future producing materially different returns. Results do not reﬂect the impact of taxes for taxable accounts or their owners. You cannot invest directly in an index. This presentation is supplemental to the composite report. Potomac claims compliance with the Global Investment Performance Standards (GIPS®). The Annual GIPS® Report is available upon request. GIPS® is a registered trademark of CFA Institute. CFA Institute does not endorse or promote this organization, nor does it warrant the accuracy or quality of the content contained herein.   '
$para5 = 'Potomac Funds are distributed by Paralel Distributors LLC. Paralel is not affiliated with Potomac Fund Management, Inc.'

$shpBody.Name = "TextBox 1"

$trBody = $shpBody.TextFrame.TextRange
$trBody.Text = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5
$trBody.Font.Size = 20
$trBody.Font.Bold = $false
$trBody.Font.Italic = $false
$trBody.Font.Name = "Quicksand (TT)"
$trBody.Font.Color.RGB = 16777215

# Target box (EMU): off (1066830, 2794040) ext (22250369, 9702789).
# Shape.Left/Top/Width/Height are single-precision points in the COM
# object model, so the literals below are the closest representable
# values to that EMU box (PowerPoint itself is subject to the same
# float32 round-trip through the OM).
$shpBody.Left = 84.00236511230469
$shpBody.Top = 220.00315856933594
$shpBody.Width = 1751.99755859375
$shpBody.Height = 763.9991455078125

# --- TextBox 2: "Disclosures" heading -----------------------------------
$shpTitle.Name = "TextBox 2"

$trTitle = $shpTitle.TextFrame.TextRange
$trTitle.Text = 'Disclosures'
$trTitle.Font.Size = 36
$trTitle.Font.Bold = $false
$trTitle.Font.Italic = $false
$trTitle.Font.Name = "Rajdhani"
$trTitle.Font.Color.RGB = 16777215

# Target box (EMU): off (1066830, 1143000) ext (5333969, 559155).
$shpTitle.Left = 84.00236511230469
$shpTitle.Top = 90.0
$shpTitle.Width = 419.9975891113281
$shpTitle.Height = 44.0279541015625
